# Applies the "Updated cryptos list" data refresh to Sheet1 (columns B-E, rows 2-51).
# D-column values that look like plain numbers are written with a leading
# apostrophe (quote-prefix) so Excel keeps them as text, matching the source data
# (e.g. "1.000" / "0.8106" must stay text, not become the number 1 / 0.8106 truncated).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.356.73"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "1.919.39"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'0.8106"
$ws.Range("E5").Value = "  +3.71%  "

$ws.Range("D6").Value = "'244.39"
$ws.Range("E6").Value = "  +1.02%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +2.50%  "

$ws.Range("D9").Value = "'27.24"
$ws.Range("E9").Value = "  +3.76%  "

$ws.Range("D10").Value = "'0.07267"
$ws.Range("E10").Value = "  +5.71%  "

$ws.Range("D11").Value = "'0.7870"
$ws.Range("E11").Value = "  +5.99%  "

$ws.Range("D12").Value = "'0.08092"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").Value = "1.928.70"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").Value = "'5.417"
$ws.Range("E14").Value = "  +4.12%  "

$ws.Range("D15").Value = "'94.74"
$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").Value = "30.368.76"
$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("E17").Value = "  +2.75%  "

$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'6.051"
$ws.Range("E18").Value = "  +2.94%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'253.79"
$ws.Range("E19").Value = "  +3.00%  "

$ws.Range("D20").Value = "'0.000007843"
$ws.Range("E20").Value = "  +1.53%  "

$ws.Range("D21").Value = "2.177.86"
$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("D23").Value = "'7.951"
$ws.Range("E23").Value = "  +16.00%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "'0.1621"
$ws.Range("E25").Value = "  +17.19%  "

$ws.Range("D26").Value = "'9.510"
$ws.Range("E26").Value = "  +2.81%  "

$ws.Range("D27").Value = "'167.65"
$ws.Range("E27").Value = "  -0.36%  "

$ws.Range("D28").Value = "'19.12"
$ws.Range("E28").Value = "  +1.33%  "

$ws.Range("E29").Value = "  +5.54%  "

$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("D31").Value = "'1.538"
$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("D32").Value = "'4.349"
$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").Value = "'4.145"
$ws.Range("E33").Value = "  +1.70%  "

$ws.Range("D34").Value = "'0.05624"
$ws.Range("E34").Value = "  +0.29%  "

$ws.Range("D35").Value = "'1.302"
$ws.Range("E35").Value = "  +3.88%  "

$ws.Range("D36").Value = "'0.7427"
$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("D38").Value = "'2.718"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("D39").Value = "'0.01956"
$ws.Range("E39").Value = "  +1.52%  "

$ws.Range("D40").Value = "'2.804"
$ws.Range("E40").Value = "  +0.66%  "

$ws.Range("D41").Value = "'0.4495"
$ws.Range("E41").Value = "  +1.80%  "

$ws.Range("D42").Value = "'73.71"
$ws.Range("E42").Value = "  +2.27%  "

$ws.Range("E43").Value = "  -2.37%  "

$ws.Range("D44").Value = "'0.8564"
$ws.Range("E44").Value = "  +1.63%  "

$ws.Range("D45").Value = "'1.937"
$ws.Range("E45").Value = "  +3.43%  "

$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("E47").Value = "  +2.81%  "

$ws.Range("D48").Value = "1.029.34"
$ws.Range("E48").Value = "  +4.07%  "

$ws.Range("E49").Value = "  +3.04%  "

$ws.Range("D50").Value = "'7.666"
$ws.Range("E50").Value = "  +1.90%  "

$ws.Range("D51").Value = "2.074.52"
$ws.Range("E51").Value = "  +0.97%  "
